$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet/tab label shown in workbook.xml "sheet name"
$ws.Name = "Through 2022-05-12"

# Update header label for the "Total" column (I1)
$ws.Range("I1").Value = "2022 (through 05-12)"

# Update May's Total column value (row 6 -> May)
$ws.Range("I6").Value = 42

# Update the Total row's Total column value (row 14 -> Total)
$ws.Range("I14").Value = 594
